$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

# Row 1 (header data row)
Set-CellText $tbl 1 1 "21÷3=7, 0"
Set-CellText $tbl 1 2 "61÷5=12, 1"
Set-CellText $tbl 1 3 "58÷7=8, 2"
Set-CellText $tbl 1 4 "63÷5=12, 3"
Set-CellText $tbl 1 5 "27÷9=3, 0"

# Row 5
Set-CellText $tbl 5 1 "97÷9=10, 7"
Set-CellText $tbl 5 2 "92÷5=18, 2"
Set-CellText $tbl 5 3 "77÷9=8, 5"
Set-CellText $tbl 5 4 "45÷9=5, 0"
Set-CellText $tbl 5 5 "92÷5=18, 2"

# Row 9
Set-CellText $tbl 9 1 "11÷2=5, 1"
Set-CellText $tbl 9 2 "61÷7=8, 5"
Set-CellText $tbl 9 3 "47÷8=5, 7"
Set-CellText $tbl 9 4 "69÷2=34, 1"
Set-CellText $tbl 9 5 "23÷5=4, 3"

# Row 13
Set-CellText $tbl 13 1 "62÷4=15, 2"
Set-CellText $tbl 13 2 "74÷8=9, 2"
Set-CellText $tbl 13 3 "36÷4=9, 0"
Set-CellText $tbl 13 4 "65÷4=16, 1"
Set-CellText $tbl 13 5 "52÷2=26, 0"

# Row 17
Set-CellText $tbl 17 1 "76÷7=10, 6"
Set-CellText $tbl 17 2 "72÷3=24, 0"
Set-CellText $tbl 17 3 "93÷6=15, 3"
Set-CellText $tbl 17 4 "83÷7=11, 6"
Set-CellText $tbl 17 5 "49÷9=5, 4"
